$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price observation was inserted right before the former row 374,
# pushing every subsequent record (old rows 374-468) down by one row.
$ws.Rows("374:374").Insert()

# Populate the newly inserted row with the new observation's data.
$ws.Range("A374").Value2 = 3
$ws.Range("B374").Value2 = "Femacal de La Calera"
$ws.Range("C374").Value2 = "Coquimbo"
$ws.Range("D374").Value2 = 44932
$ws.Range("E374").Value2 = 5
$ws.Range("F374").Value2 = 100112009
$ws.Range("G374").Value2 = "Acelga"
$ws.Range("H374").Value2 = "Sin especificar"
$ws.Range("I374").Value2 = "Primera"
$ws.Range("J374").Value2 = 200
$ws.Range("K374").Value2 = 3500
$ws.Range("L374").Value2 = 3800
$ws.Range("M374").Value2 = 3635
$ws.Range("N374").Value2 = "$/docena de atados (6 kilos)"
$ws.Range("O374").Value2 = "Provincia de Quillota"
$ws.Range("P374").Value2 = 606
$ws.Range("Q374").Value2 = 6
$ws.Range("R374").Value2 = "Hortaliza"
